$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")

# Update the External Contact test data row
$ws.Range("A2").Value = "Activity Test External Contact"
$ws.Range("B2").Value = "ActivityCompany"

# Make Contact the active sheet/tab with A2:B2 selected
$ws.Activate()
$ws.Range("A2:B2").Select()
